# Update the "Data Warehouse" documentation sheet:
#  - DimAgent.HireDate  -> DimAgent.WorkingDuration (datetime -> int)
#  - DimAgent.BirthDate -> DimAgent.Age              (datetime -> int)
#  - Fill in the previously-blank Destination Type / Dimension Column
#    cells for the FactDaysOnMarket.LocationKey row.
#  - Move the saved selection/scroll position to D34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Warehouse")

# Row 10: FactDaysOnMarket.LocationKey - fill in D (Dimension Column) and E (Destination Type)
$ws.Range("D10").Value = "na"
$ws.Range("E10").Value = "int"

# Row 29: HireDate -> WorkingDuration, type datetime -> int
$ws.Range("A29").Value = "DWRedwood.dbo.DimAgent.WorkingDuration"
$ws.Range("E29").Value = "int"

# Row 30: BirthDate -> Age, type datetime -> int
$ws.Range("A30").Value = "DWRedwood.dbo.DimAgent.Age"
$ws.Range("E30").Value = "int"

# Update view: scroll/selection moved to D34 (row 13 at top of viewport)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
